$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new row (row 4) that duplicates the "sauceDemoLoginTest" test case
# already present in row 2 (testcasename / username / password).
$ws.Range("A4").Value = "sauceDemoLoginTest"
$ws.Range("B4").Value = "standard_user"
$ws.Range("C4").Value = "secret_sauce"

# Update the view: zoom in to 145% and move the active cell/selection.
$ws.Select()
$ws.Range("C7").Select()
$excel.ActiveWindow.Zoom = 145
